$d = $word.ActiveDocument

$startPos = $d.Paragraphs.Item(66).Range.Start
$endPos = $d.Paragraphs.Item(70).Range.End
$r = $d.Range($startPos, $endPos)

$body = @'
<w:p><w:pPr>
        <w:spacing w:line="480" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">2016-4-19 I updated the input file to make everything combines into one. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">I also add lots of running controls in the model input file. Note calculation </w:t>
      </w:r>
      <w:r>
        <w:t>expressions (non-numbers)</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">in the excel file like c = a* km, </w:t>
      </w:r>
      <w:r>
        <w:t>are not readable</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> by the model.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="480" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">2016-4-21 I finished the sensitivity file of the model. The parameter values were written in an array in the script, running conditions, like temperature, RH, radiation and </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>swc</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> were combined in a excel file that read by the model. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="480" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">The model output table were also updated, and moved into the folder of model scenarios. Next time, I will not list the information of all the leaves, I will just output the leaf numbers that we are </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>interested,</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> else the table will be too big. R code is available for converting the excel format into the text format for inputting into the model. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Do not change the existing sequence anymore. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="480" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="480" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">2016-4-26 </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>Add</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> a PLC (percentage loss of conductivity) function into the water flux optimization part to limit the transpiration under low water potential condition. The idea is discussed with Greg. We found it is necessary to decrease </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>gs_min</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> in order to reach zero water flux under limiting condition. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="480" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">2016-4-28 as it is really difficult to optimize transpiration, water potential and ABA simultaneously, due to different time scales in those data and the method. An approximation is taken by optimizing the leaf water potential during the midday. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="480" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">2016-5-31 I corrected g0 in the calculation of </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>gs</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> by the Tardieu’s method. As in Tardieu’s method, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>gs</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>is</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> the conductance for water vapor, and the g0 initially is for CO2. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:line="480" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t>I also add the percentage of loss of conductivity</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> (PLC)</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">as a function of leaf water potential </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">in </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">the calculation of </w:t>
      </w:r>
      <w:r>
        <w:t>leaf conductance</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">. There is interdependence between leaf water potential and leaf conductance, </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>one iteration</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> is used in the calculation. Leaf conductance was first calculated by the water flux and the water potential of the previous step. Leaf water potential was updated by this leaf conductance. Then leaf conductance and leaf water potential was calculated again. </w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
'@

$xmlFrag = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>" + $body + "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

$r.InsertXML($xmlFrag)

Write-Output "Done applying edit"
